$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) per the commit diff.
# D/E columns are plain numeric-looking text (e.g. "8.20", "0.0000269") stored as
# inline strings in the source sheet, so NumberFormat is forced to "@" (Text) before
# the assignment and reset to "General" after, to stop Excel from auto-coercing the
# text into a float and silently dropping trailing zeros / thousands-style dots.
$updates = @(
    @{ Cell = "D2"; Value = '67.032.47' }
    @{ Cell = "E2"; Value = '  +1.19%  ' }
    @{ Cell = "D3"; Value = '3.273.57' }
    @{ Cell = "E3"; Value = '  -2.03%  ' }
    @{ Cell = "D4"; Value = '0.998' }
    @{ Cell = "E4"; Value = '  -0.22%  ' }
    @{ Cell = "D5"; Value = '575.68' }
    @{ Cell = "E5"; Value = '  -1.70%  ' }
    @{ Cell = "D6"; Value = '171.11' }
    @{ Cell = "E6"; Value = '  -7.77%  ' }
    @{ Cell = "E7"; Value = '  -0.02%  ' }
    @{ Cell = "D8"; Value = '0.576' }
    @{ Cell = "E8"; Value = '  +0.04%  ' }
    @{ Cell = "D9"; Value = '3.261.65' }
    @{ Cell = "E9"; Value = '  -2.31%  ' }
    @{ Cell = "D10"; Value = '0.172' }
    @{ Cell = "E10"; Value = '  -5.03%  ' }
    @{ Cell = "D11"; Value = '0.567' }
    @{ Cell = "E11"; Value = '  -2.75%  ' }
    @{ Cell = "D12"; Value = '44.61' }
    @{ Cell = "E12"; Value = '  -4.98%  ' }
    @{ Cell = "D13"; Value = '0.0000269' }
    @{ Cell = "E13"; Value = '  -0.01%  ' }
    @{ Cell = "D14"; Value = '685.37' }
    @{ Cell = "E14"; Value = '  +2.52%  ' }
    @{ Cell = "D15"; Value = '3.786.00' }
    @{ Cell = "E15"; Value = '  -2.35%  ' }
    @{ Cell = "D16"; Value = '8.20' }
    @{ Cell = "E16"; Value = '  -3.68%  ' }
    @{ Cell = "D17"; Value = '66.891.72' }
    @{ Cell = "E17"; Value = '  +0.74%  ' }
    @{ Cell = "E18"; Value = '  +0.42%  ' }
    @{ Cell = "D19"; Value = '3.259.29' }
    @{ Cell = "E19"; Value = '  -2.51%  ' }
    @{ Cell = "D20"; Value = '17.11' }
    @{ Cell = "E20"; Value = '  -4.32%  ' }
    @{ Cell = "D21"; Value = '10.62' }
    @{ Cell = "E21"; Value = '  -4.25%  ' }
    @{ Cell = "D22"; Value = '0.877' }
    @{ Cell = "E22"; Value = '  -2.17%  ' }
    @{ Cell = "D23"; Value = '16.80' }
    @{ Cell = "E23"; Value = '  -4.96%  ' }
    @{ Cell = "D24"; Value = '5.18' }
    @{ Cell = "E24"; Value = '  +2.67%  ' }
    @{ Cell = "D25"; Value = '97.68' }
    @{ Cell = "E25"; Value = '  -3.55%  ' }
    @{ Cell = "E26"; Value = '  -4.60%  ' }
    @{ Cell = "D27"; Value = '2.62' }
    @{ Cell = "E27"; Value = '  -5.76%  ' }
    @{ Cell = "D28"; Value = '33.11' }
    @{ Cell = "E28"; Value = '  +2.63%  ' }
    @{ Cell = "D29"; Value = '9.00' }
    @{ Cell = "E29"; Value = '  -4.53%  ' }
    @{ Cell = "D30"; Value = '8.25' }
    @{ Cell = "E30"; Value = '  -3.04%  ' }
    @{ Cell = "D31"; Value = '6.55' }
    @{ Cell = "E31"; Value = '  -3.47%  ' }
    @{ Cell = "D32"; Value = '575.70' }
    @{ Cell = "E32"; Value = '  -6.44%  ' }
    @{ Cell = "D33"; Value = '10.74' }
    @{ Cell = "E33"; Value = '  -3.23%  ' }
    @{ Cell = "D34"; Value = '3.810.52' }
    @{ Cell = "E34"; Value = '  -1.26%  ' }
    @{ Cell = "B35"; Value = 'Hedera' }
    @{ Cell = "C35"; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = "D35"; Value = '0.102' }
    @{ Cell = "E35"; Value = '  -3.40%  ' }
    @{ Cell = "B36"; Value = 'Dai' }
    @{ Cell = "C36"; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = "D36"; Value = '0.999' }
    @{ Cell = "E36"; Value = '  -0.20%  ' }
    @{ Cell = "D37"; Value = '55.28' }
    @{ Cell = "E37"; Value = '  -1.70%  ' }
    @{ Cell = "D38"; Value = '3.26' }
    @{ Cell = "E38"; Value = '  -16.12%  ' }
    @{ Cell = "D39"; Value = '0.127' }
    @{ Cell = "E39"; Value = '  -1.07%  ' }
    @{ Cell = "D40"; Value = '3.37' }
    @{ Cell = "E40"; Value = '  -1.22%  ' }
    @{ Cell = "D41"; Value = '2.55' }
    @{ Cell = "E41"; Value = '  -4.37%  ' }
    @{ Cell = "D42"; Value = '31.20' }
    @{ Cell = "E42"; Value = '  -4.95%  ' }
    @{ Cell = "D43"; Value = '0.0₃0653' }
    @{ Cell = "E43"; Value = '  -6.83%  ' }
    @{ Cell = "D44"; Value = '2.95' }
    @{ Cell = "E44"; Value = '  -7.42%  ' }
    @{ Cell = "D45"; Value = '0.322' }
    @{ Cell = "E45"; Value = '  -4.50%  ' }
    @{ Cell = "D46"; Value = '0.0398' }
    @{ Cell = "E46"; Value = '  -4.69%  ' }
    @{ Cell = "E47"; Value = '  -0.12%  ' }
    @{ Cell = "D48"; Value = '0.126' }
    @{ Cell = "E48"; Value = '  -1.73%  ' }
    @{ Cell = "D49"; Value = '2.51' }
    @{ Cell = "E49"; Value = '  -1.52%  ' }
    @{ Cell = "E50"; Value = '  +3.10%  ' }
    @{ Cell = "D51"; Value = '127.65' }
    @{ Cell = "E51"; Value = '  -1.27%  ' }
)

foreach ($u in $updates) {
    $col = $u.Cell -replace '[0-9]+$', ''
    $range = $ws.Range($u.Cell)
    if ($col -eq "D" -or $col -eq "E") {
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.NumberFormat = "General"
    } else {
        $range.Value = $u.Value
    }
}
